# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which contain identical data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2 = 1390
    3 = 2304
    4 = 399
    6 = 6461
    7 = 322
    8 = 126
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
